$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix existing rows (remove "loss" column D1, fix row 31 values) ---
$ws.Range("D1").ClearContents()
$ws.Range("B31").Value = 30
$ws.Range("E31").Value = 6000

# --- Step 2: add new rows, in chronological order for shared-string table fidelity ---
# 2a: write text/label cells in original authoring order (affects shared string table order)
$ws.Range("H51").Value = "training with input optim as well "
$ws.Range("A51").Value = "adversarial + input optim"
$ws.Range("A54").Value = "denoising "
$ws.Range("A55").Value = "snail"
$ws.Range("A64").Value = "results "
$ws.Range("A66").Value = "kate - regular"
$ws.Range("A67").Value = "input - optim"
$ws.Range("A68").Value = "adversarial"
$ws.Range("A70").Value = "barbara - regular"
$ws.Range("A71").Value = "input - optim"
$ws.Range("A72").Value = "adversarial"
$ws.Range("A75").Value = "kate - regular"
$ws.Range("A76").Value = "input - optim"
$ws.Range("A77").Value = "adversarial"
$ws.Range("A79").Value = "vase - regular"
$ws.Range("A80").Value = "input - optim"
$ws.Range("A81").Value = "adversarial"
$ws.Range("A83").Value = "library - regular"
$ws.Range("A84").Value = "input - optim"
$ws.Range("A85").Value = "adversarial"
$ws.Range("A87").Value = "denoising"
$ws.Range("A87").Font.Bold = $true
$ws.Range("A88").Value = "snail - regular"
$ws.Range("A89").Value = "input - optim"
$ws.Range("A90").Value = "adversarial"
$ws.Range("A92").Value = "jet - regular"
$ws.Range("A93").Value = "input - optim"
$ws.Range("A94").Value = "adversarial"
$ws.Range("A96").Value = "sr"
$ws.Range("A96").Font.Bold = $true
$ws.Range("A97").Value = "zebra factor 4 - regular"
$ws.Range("A98").Value = "input - optim"
$ws.Range("A99").Value = "adversarial"
$ws.Range("A101").Value = "zebra factor 8 - regular"
$ws.Range("A102").Value = "input - optim"
$ws.Range("A103").Value = "adversarial"
$ws.Range("C64").Value = "psnr_masked"
$ws.Range("D64").Value = "num_iter"
$ws.Range("E93").Value = "worst looking results"

# 2b: write remaining numeric cells for each new row
# Row 51
$ws.Range("B51").Value = 40.82
$ws.Range("B51").Font.Bold = $true
$ws.Range("C51").Value = 43.0
$ws.Range("E51").Value = 19400.0
$ws.Range("F51").Value = 0.001
$ws.Range("G51").Value = 0.0001
# Row 54
# Row 55
$ws.Range("B55").Value = 20.52
$ws.Range("E55").Value = 8250.0
$ws.Range("F55").Value = 0.001
$ws.Range("G55").Value = 0.005
# Row 64
$ws.Range("B64").Value = "psnr"
# Row 65
$ws.Range("A65").Value = "restoration"
$ws.Range("A65").Font.Bold = $true
# Row 66
$ws.Range("B66").Value = 24.75
$ws.Range("C66").Value = 48.44
$ws.Range("D66").Value = 1000.0
# Row 67
$ws.Range("B67").Value = 25.01
$ws.Range("C67").Value = 48.77
$ws.Range("D67").Value = 1000.0
# Row 68
$ws.Range("B68").Value = 24.04
$ws.Range("C68").Value = 58.8
$ws.Range("D68").Value = 4600.0
# Row 70
$ws.Range("B70").Value = 32.12
$ws.Range("C70").Value = 41.23
$ws.Range("D70").Value = 11000.0
# Row 71
$ws.Range("B71").Value = 29.66
$ws.Range("C71").Value = 44.75
$ws.Range("D71").Value = 11000.0
# Row 72
$ws.Range("B72").Value = 28.3
$ws.Range("C72").Value = 40.12
$ws.Range("D72").Value = 22000.0
# Row 74
$ws.Range("A74").Value = "inpainting"
$ws.Range("A74").Font.Bold = $true
# Row 75
$ws.Range("B75").Value = 38.74
$ws.Range("C75").Value = 39.89
$ws.Range("D75").Value = 6000.0
# Row 76
$ws.Range("B76").Value = 42.27
$ws.Range("C76").Value = 49.28
$ws.Range("D76").Value = 6000.0
# Row 77
$ws.Range("B77").Value = 39.75
$ws.Range("C77").Value = 42.5
$ws.Range("D77").Value = 15000.0
# Row 79
$ws.Range("B79").Value = 29.22
$ws.Range("C79").Value = 35.99
$ws.Range("D79").Value = 5000.0
# Row 80
$ws.Range("B80").Value = 28.86
$ws.Range("C80").Value = 35.7
$ws.Range("D80").Value = 5000.0
# Row 81
$ws.Range("B81").Value = 30.0
$ws.Range("C81").Value = 41.0
$ws.Range("D81").Value = 5000.0
# Row 83
$ws.Range("B83").Value = 19.43
$ws.Range("C83").Value = 29.53
$ws.Range("D83").Value = 3000.0
# Row 84
$ws.Range("B84").Value = 19.17
$ws.Range("C84").Value = 28.85
$ws.Range("D84").Value = 3000.0
# Row 85
$ws.Range("B85").Value = 16.8
$ws.Range("C85").Value = 23.0
$ws.Range("D85").Value = 6000.0
# Row 87
# Row 88
$ws.Range("B88").Value = 26.51
$ws.Range("D88").Value = 2400.0
# Row 89
$ws.Range("B89").Value = 29.01
$ws.Range("D89").Value = 2400.0
# Row 90
$ws.Range("B90").Value = 18.21
$ws.Range("D90").Value = 3000.0
# Row 92
$ws.Range("B92").Value = 31.6
$ws.Range("C92").Value = 27.54
$ws.Range("D92").Value = 2400.0
# Row 93
$ws.Range("B93").Value = 35.23
$ws.Range("C93").Value = 29.31
$ws.Range("D93").Value = 2400.0
# Row 94
$ws.Range("B94").Value = 19.0
$ws.Range("C94").Value = 18.87
$ws.Range("D94").Value = 2200.0
# Row 96
# Row 97
$ws.Range("B97").Value = 24.03
$ws.Range("C97").Value = 32.83
$ws.Range("D97").Value = 2000.0
# Row 98
$ws.Range("B98").Value = 23.0
$ws.Range("C98").Value = 43.0
$ws.Range("D98").Value = 2000.0
# Row 99
$ws.Range("D99").Value = 1800.0
# Row 101
$ws.Range("B101").Value = 19.54
$ws.Range("C101").Value = 42.12
$ws.Range("D101").Value = 4000.0
# Row 102
$ws.Range("B102").Value = 19.12
$ws.Range("C102").Value = 57.19
$ws.Range("D102").Value = 4000.0
# Row 103
$ws.Range("B103").Value = 17.93
$ws.Range("C103").Value = 35.56
$ws.Range("D103").Value = 4000.0

# --- Step 3: column width for column E ---
$ws.Columns.Item(5).ColumnWidth = 17

# --- Step 4: view settings ---
$win = $excel.ActiveWindow
$win.ScrollRow = 81
$win.ScrollColumn = 1
$ws.Range("C100").Select()
